$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# ALC
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H33").Value = 840.69696
$ws.Range("I33").Value = 578.4138
$ws.Range("J33").Value = 2742.25
$ws.Range("K33").Value = 578.4138
$ws.Range("L33").Value = 2742.25
$ws.Range("M33").Value = -349.4138
$ws.Range("N33").Value = -3200.25

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H70").Value = 400
$ws.Range("I70").Value = 400
$ws.Range("K70").Value = 1200
$ws.Range("M70").Value = -930

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H73").Value = 400
$ws.Range("I73").Value = 400
$ws.Range("K73").Value = 1200
$ws.Range("M73").Value = -264

$ws.Range("H107").Value = 35714616
$ws.Range("I107").Value = 50000230
$ws.Range("J107").Value = 568.25
$ws.Range("K107").Value = 50000230
$ws.Range("L107").Value = 568.25
$ws.Range("M107").Value = -49998310
$ws.Range("N107").Value = -4408.25

$ws.Range("H125").Value = 2916.6667
$ws.Range("J125").Value = 3375
$ws.Range("L125").Value = 30375
$ws.Range("N125").Value = -35295

$ws.Range("H137").Value = 2158.558
$ws.Range("I137").Value = 2124.4348
$ws.Range("J137").Value = 2197.8
$ws.Range("K137").Value = 6373.3044
$ws.Range("L137").Value = 6593.400000000001
$ws.Range("M137").Value = -3823.3044
$ws.Range("N137").Value = -11693.4

$ws.Range("H138").Value = 2524.541
$ws.Range("I138").Value = 1069.4117
$ws.Range("J138").Value = 3086.75
$ws.Range("K138").Value = 3208.2351
$ws.Range("L138").Value = 9260.25
$ws.Range("M138").Value = 1931.7649
$ws.Range("N138").Value = -19540.25

$ws.Range("H141").Value = 21801704
$ws.Range("I141").Value = 3336445.5
$ws.Range("J141").Value = 83352570
$ws.Range("K141").Value = 10009336.5
$ws.Range("L141").Value = 250057710
$ws.Range("M141").Value = -10004156.5
$ws.Range("N141").Value = -250068070

# ------------------------------------------------------------------
# ARM
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 1266.5
$ws.Range("I2").Value = 1319.8
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1319.8
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -1206.8
$ws.Range("N2").Value = -1226

$ws.Range("H32").Value = 8731.475
$ws.Range("I32").Value = 3071.4075
$ws.Range("J32").Value = 20487
$ws.Range("K32").Value = 3071.4075
$ws.Range("L32").Value = 20487
$ws.Range("M32").Value = -2784.4075
$ws.Range("N32").Value = -21061

$ws.Range("H116").Value = 1266.5
$ws.Range("I116").Value = 1319.8
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 1319.8
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 974.2
$ws.Range("N116").Value = -5588

$ws.Range("H122").Value = 4006.6
$ws.Range("I122").Value = 3692.8
$ws.Range("J122").Value = 4948
$ws.Range("K122").Value = 11078.4
$ws.Range("L122").Value = 14844
$ws.Range("M122").Value = -8628.400000000001
$ws.Range("N122").Value = -19744

# ------------------------------------------------------------------
# BSM
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 1266.5
$ws.Range("I3").Value = 1319.8
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 1319.8
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -1205.8
$ws.Range("N3").Value = -1228

$ws.Range("H80").Value = 132.4
$ws.Range("I80").Value = 142.4
$ws.Range("K80").Value = 142.4
$ws.Range("M80").Value = 855.6

$ws.Range("H83").Value = 132.4
$ws.Range("I83").Value = 142.4
$ws.Range("K83").Value = 712
$ws.Range("M83").Value = 4280

$ws.Range("H86").Value = 40004160
$ws.Range("I86").Value = 50002700
$ws.Range("K86").Value = 50002700
$ws.Range("M86").Value = -50001577

$ws.Range("H89").Value = 40004160
$ws.Range("I89").Value = 50002700
$ws.Range("K89").Value = 250013500
$ws.Range("M89").Value = -250007884

# ------------------------------------------------------------------
# CRP
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 45456390
$ws.Range("I16").Value = 71429990
$ws.Range("J16").Value = 2599.75
$ws.Range("K16").Value = 71429990
$ws.Range("L16").Value = 2599.75
$ws.Range("M16").Value = -71429703
$ws.Range("N16").Value = -3173.75

$ws.Range("H113").Value = 45456390
$ws.Range("I113").Value = 71429990
$ws.Range("J113").Value = 2599.75
$ws.Range("K113").Value = 71429990
$ws.Range("L113").Value = 2599.75
$ws.Range("M113").Value = -71427820
$ws.Range("N113").Value = -6939.75

# ------------------------------------------------------------------
# CUL
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H31").Value = 4999
$ws.Range("J31").Value = 4999
$ws.Range("L31").Value = 14997
$ws.Range("N31").Value = -15573

$ws.Range("H39").Value = 2817.2727
$ws.Range("J39").Value = 2817.2727
$ws.Range("L39").Value = 8451.8181
$ws.Range("N39").Value = -9039.8181

$ws.Range("H131").Value = 1235453.1
$ws.Range("I131").Value = 3704384.8
$ws.Range("J131").Value = 987.3611
$ws.Range("K131").Value = 11113154.4
$ws.Range("L131").Value = 2962.0833
$ws.Range("M131").Value = -11108114.4
$ws.Range("N131").Value = -13042.0833

# ------------------------------------------------------------------
# GSM
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H80").Value = 3322.4138
$ws.Range("I80").Value = 3724.6667
$ws.Range("J80").Value = 2891.4285
$ws.Range("K80").Value = 3724.6667
$ws.Range("L80").Value = 2891.4285
$ws.Range("M80").Value = -2726.6667
$ws.Range("N80").Value = -4887.4285

$ws.Range("H83").Value = 3322.4138
$ws.Range("I83").Value = 3724.6667
$ws.Range("J83").Value = 2891.4285
$ws.Range("K83").Value = 18623.3335
$ws.Range("L83").Value = 14457.1425
$ws.Range("M83").Value = -13631.3335
$ws.Range("N83").Value = -24441.1425

$ws.Range("H102").Value = 1260.7142
$ws.Range("I102").Value = 1287.4546
$ws.Range("J102").Value = 1162.6666
$ws.Range("K102").Value = 1287.4546
$ws.Range("L102").Value = 1162.6666
$ws.Range("M102").Value = 334.5454
$ws.Range("N102").Value = -4406.6666

$ws.Range("H126").Value = 2478.92
$ws.Range("I126").Value = 2616.1765
$ws.Range("J126").Value = 2187.25
$ws.Range("K126").Value = 7848.529500000001
$ws.Range("L126").Value = 6561.75
$ws.Range("M126").Value = -5378.529500000001
$ws.Range("N126").Value = -11501.75

# ------------------------------------------------------------------
# LTW
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H46").Value = 1211.9131
$ws.Range("I46").Value = 1127.8334
$ws.Range("J46").Value = 1303.6364
$ws.Range("K46").Value = 1127.8334
$ws.Range("L46").Value = 1303.6364
$ws.Range("M46").Value = -939.8334
$ws.Range("N46").Value = -1679.6364

$ws.Range("H82").Value = 10103624
$ws.Range("I82").Value = 18182758
$ws.Range("J82").Value = 4706.5
$ws.Range("K82").Value = 18182758
$ws.Range("L82").Value = 4706.5
$ws.Range("M82").Value = -18182397
$ws.Range("N82").Value = -5428.5

$ws.Range("H85").Value = 10103624
$ws.Range("I85").Value = 18182758
$ws.Range("J85").Value = 4706.5
$ws.Range("K85").Value = 18182758
$ws.Range("L85").Value = 4706.5
$ws.Range("M85").Value = -18181510
$ws.Range("N85").Value = -7202.5

$ws.Range("H132").Value = 1846797
$ws.Range("I132").Value = 3690285.5
$ws.Range("J132").Value = 3308.4119
$ws.Range("K132").Value = 11070856.5
$ws.Range("L132").Value = 9925.235700000001
$ws.Range("M132").Value = -11068326.5
$ws.Range("N132").Value = -14985.2357

# ------------------------------------------------------------------
# WVR
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H138").Value = 20000
$ws.Range("J138").Value = 20000
$ws.Range("L138").Value = 20000
$ws.Range("N138").Value = -30280
